# Apply updated cryptocurrency price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.703.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.251.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.48%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -4.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.245.70"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.44%  "
$ws.Range("E10").Value = "  -8.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.569"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000268"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "660.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.773.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.683.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.243.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.875"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.20%  "
$ws.Range("E26").Value = "  -5.06%  "
$ws.Range("E27").Value = "  -7.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "31.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.13%  "
$ws.Range("E31").Value = "  -2.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "566.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.749.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.68%  "
$ws.Range("E35").Value = "  -4.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -14.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.129"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "32.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.76%  "
$ws.Range("E41").Value = "  -8.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0655"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.21%  "
$ws.Range("E43").Value = "  -8.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.324"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.97%  "
$ws.Range("E46").Value = "  -6.88%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("E48").Value = "  -3.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.95%  "
